# Add a new bullet item "User can add a row using Post insert row API"
# right after the "Using Get Table API ..." bullet, matching the
# surrounding ListParagraph / numbered-list formatting.

$d = $word.ActiveDocument

# Locate the paragraph whose text is the "Using Get Table API..." bullet.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Using Get Table API user can see user specific table data*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Using Get Table API...' paragraph"
}

# Insert a new paragraph right after it; InsertParagraphAfter clones the
# paragraph formatting (style, numbering, rPr) of the paragraph it is
# called on, so the new bullet keeps the same list style/formatting.
$target.Range.InsertParagraphAfter()

# Find the paragraph that was just created (the one now following target)
$newPara = $target.Next()
$newPara.Range.Text = "User can add a row using Post insert row API"

Write-Output "Inserted new bullet after 'Using Get Table API...' paragraph"
